# Update the NATMI ligand-receptor pair data (Il17c-Il17re) with the new TPM-based
# values. The sending/target-cluster matrix grows from 3 targets to the full 4 targets
# per sending cluster (16 rows total), and all computed metrics are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Il17c', 'Il17re', 'ECs', 3, 1, 0.8238633333333333, 2.47159, 0.2181573957783601, 0.2181573957783601, 3, 1, 0.09324366666666667, 0.279731, 0.1594970327830716, 0.1594970327830716, 0.07682003803222222, 0.69138034229, 0.03479545730633063, 0.03479545730633063),
    @('ECs', 'Il17c', 'Il17re', 'FAPs', 3, 1, 0.8238633333333333, 2.47159, 0.2181573957783601, 0.2181573957783601, 3, 1, 0.261563, 0.784689, 0.4474140054463598, 0.4474140054463597, 0.2154921650566667, 1.93942948551, 0.09760667426294288, 0.09760667426294288),
    @('ECs', 'Il17c', 'Il17re', 'MuSCs', 3, 1, 0.8238633333333333, 2.47159, 0.2181573957783601, 0.2181573957783601, 3, 1, 0.142436, 0.427308, 0.2436424925534487, 0.2436424925534487, 0.1173477977466667, 1.05613017972, 0.05315241167640887, 0.05315241167640887),
    @('ECs', 'Il17c', 'Il17re', 'Resolving-Mac', 3, 1, 0.8238633333333333, 2.47159, 0.2181573957783601, 0.2181573957783601, 2, 0.6666666666666666, 0.087368, 0.262104, 0.14944646921712, 0.14944646921712, 0.07197929170666666, 0.64781362536, 0.03260285253267776, 0.03260285253267776),
    @('FAPs', 'Il17c', 'Il17re', 'ECs', 3, 1, 1.210478333333333, 3.631435, 0.3205322899584435, 0.3205322899584435, 3, 1, 0.09324366666666667, 0.279731, 0.1594970327830716, 0.1594970327830716, 0.1128694382205556, 1.015824943985, 0.05112394915953487, 0.05112394915953487),
    @('FAPs', 'Il17c', 'Il17re', 'FAPs', 3, 1, 1.210478333333333, 3.631435, 0.3205322899584435, 0.3205322899584435, 3, 1, 0.261563, 0.784689, 0.4474140054463598, 0.4474140054463597, 0.3166163443016666, 2.849547098715, 0.1434106357252012, 0.1434106357252012),
    @('FAPs', 'Il17c', 'Il17re', 'MuSCs', 3, 1, 1.210478333333333, 3.631435, 0.3205322899584435, 0.3205322899584435, 3, 1, 0.142436, 0.427308, 0.2436424925534487, 0.2436424925534487, 0.1724156918866667, 1.55174122698, 0.07809528606933992, 0.07809528606933992),
    @('FAPs', 'Il17c', 'Il17re', 'Resolving-Mac', 3, 1, 1.210478333333333, 3.631435, 0.3205322899584435, 0.3205322899584435, 2, 0.6666666666666666, 0.087368, 0.262104, 0.14944646921712, 0.14944646921712, 0.1057570710266667, 0.9518136392400001, 0.04790241900436751, 0.0479024190043675),
    @('MuSCs', 'Il17c', 'Il17re', 'ECs', 3, 1, 0.316433, 0.9492990000000001, 0.0837908381466997, 0.0837908381466997, 3, 1, 0.09324366666666667, 0.279731, 0.1594970327830716, 0.1594970327830716, 0.02950537317433334, 0.2655483585690001, 0.01336439005880521, 0.01336439005880521),
    @('MuSCs', 'Il17c', 'Il17re', 'FAPs', 3, 1, 0.316433, 0.9492990000000001, 0.0837908381466997, 0.0837908381466997, 3, 1, 0.261563, 0.784689, 0.4474140054463598, 0.4474140054463597, 0.08276716477900001, 0.7449044830110001, 0.03748919451492255, 0.03748919451492255),
    @('MuSCs', 'Il17c', 'Il17re', 'MuSCs', 3, 1, 0.316433, 0.9492990000000001, 0.0837908381466997, 0.0837908381466997, 3, 1, 0.142436, 0.427308, 0.2436424925534487, 0.2436424925534487, 0.045071450788, 0.405643057092, 0.02041500865920451, 0.02041500865920451),
    @('MuSCs', 'Il17c', 'Il17re', 'Resolving-Mac', 3, 1, 0.316433, 0.9492990000000001, 0.0837908381466997, 0.0837908381466997, 2, 0.6666666666666666, 0.087368, 0.262104, 0.14944646921712, 0.14944646921712, 0.027646118344, 0.248815065096, 0.01252224491376744, 0.01252224491376744),
    @('Resolving-Mac', 'Il17c', 'Il17re', 'ECs', 3, 1, 1.425688333333333, 4.277065, 0.3775194761164967, 0.3775194761164967, 3, 1, 0.09324366666666667, 0.279731, 0.1594970327830716, 0.1594970327830716, 0.1329364077238889, 1.196427669515, 0.06021323625840088, 0.06021323625840088),
    @('Resolving-Mac', 'Il17c', 'Il17re', 'FAPs', 3, 1, 1.425688333333333, 4.277065, 0.3775194761164967, 0.3775194761164967, 3, 1, 0.261563, 0.784689, 0.4474140054463598, 0.4474140054463597, 0.3729073175316667, 3.356165857785, 0.1689075009432931, 0.1689075009432931),
    @('Resolving-Mac', 'Il17c', 'Il17re', 'MuSCs', 3, 1, 1.425688333333333, 4.277065, 0.3775194761164967, 0.3775194761164967, 3, 1, 0.142436, 0.427308, 0.2436424925534487, 0.2436424925534487, 0.2030693434466667, 1.82762409102, 0.0919797861484954, 0.09197978614849539),
    @('Resolving-Mac', 'Il17c', 'Il17re', 'Resolving-Mac', 3, 1, 1.425688333333333, 4.277065, 0.3775194761164967, 0.3775194761164967, 2, 0.6666666666666666, 0.087368, 0.262104, 0.14944646921712, 0.14944646921712, 0.1245595383066667, 1.12103584476, 0.0564189527663073, 0.05641895276630728)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
